# StructureDefinition-ror-organization-financial-help-type.xlsx update
# - Metadata sheet: Version/Status/Date/Contact refresh for 0.4.0-snapshot-1
# - Elements sheet: swap the two "Mapping" columns (AK <-> AL), header + data

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Metadata sheet
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value  = "0.4.0-snapshot-1"
$meta.Range("B6").Value  = "draft"
$meta.Range("B8").Value  = "2024-05-23T12:16:26+00:00"
$meta.Range("B10").Value = "ANS (https://esante.gouv.fr)"

# ---------------------------------------------------------------------------
# Elements sheet - swap columns AK (37) and AL (38): the "Mapping: RIM
# Mapping" column and the "Mapping: Spécification métier vers l'extension
# ROR FinancialHelpType" column trade places, header and data alike.
# ---------------------------------------------------------------------------
$el = $wb.Worksheets.Item("Elements")

# Row 1 - header
$el.Range("AK1").Value = "Mapping: Spécification métier vers l'extension ROR FinancialHelpType"
$el.Range("AL1").Value = "Mapping: RIM Mapping"

# Row 2 - both empty already, nothing to change

# Row 3
$el.Range("AK3").ClearContents()
$el.Range("AL3").Value = "n/a"

# Row 4 - both empty already, nothing to change

# Row 5
$el.Range("AK5").ClearContents()
$el.Range("AL5").Value = "N/A"

# Row 6
$el.Range("AK6").Value = "aideFinanciere"
$el.Range("AL6").Value = "N/A"

# ---------------------------------------------------------------------------
# Column widths follow the swapped content (AK now holds the long French
# label, AL now holds the short RIM mapping values).
# ---------------------------------------------------------------------------
$el.Columns.Item(37).ColumnWidth = 74.0703125
$el.Columns.Item(38).ColumnWidth = 24.98046875
